$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# --- Data changes ---

# End Year 2090 -> 2069
$ws.Range("B4").Value = 2069

# fix_demand_to_representative_year: TRUE -> FALSE
$ws.Range("B28").Value = $false

# fix_profiles_to_representative_year: TRUE -> FALSE
$ws.Range("B29").Value = $false

# last_years_IRR_or_NPV (B31): TRUE -> FALSE
$ws.Range("B31").Value = $false

# capacity_remuneration_mechanism: none -> strategic_reserve_ger
$ws.Range("B44").Value = "strategic_reserve_ger"

# New check row 48: demand/weather-year correlation check
$ws.Range("B48").Formula = '=IF(OR(AND(B29=TRUE,B28=FALSE),AND(B29=FALSE,B28=TRUE)),"demand must be correlated with weather year","ok")'

# --- Conditional formatting: move the "ok" check range from B49:B53 to B48:B53 ---
$oldRule = $ws.Range("B49:B53").FormatConditions.Item(1)
$oldRule.Delete()
$newRule = $ws.Range("B48:B53").FormatConditions.Add(1, 4, '"ok"')
$newRule.Interior.Color = 255

# --- View changes: scroll so row 14 is at top, select C35 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("C35").Select()
